$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.142.68"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "'1.919.78"
$ws.Range("E3").Value = "  +2.63%  "

$ws.Range("D5").Value = "'319.58"
$ws.Range("E5").Value = "  +0.06%  "

$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "'0.5066"
$ws.Range("E7").Value = "  -0.31%  "

$ws.Range("D8").Value = "'0.4071"
$ws.Range("E8").Value = "  +3.38%  "

$ws.Range("D9").Value = "'0.08350"
$ws.Range("E9").Value = "  +2.17%  "

$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.111"
$ws.Range("E10").Value = "  +1.70%  "

$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'42.33"
$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("D12").Value = "'24.04"
$ws.Range("E12").Value = "  +5.87%  "

$ws.Range("D13").Value = "'6.432"
$ws.Range("E13").Value = "  +2.64%  "

$ws.Range("D14").Value = "'1.916.30"
$ws.Range("E14").Value = "  +2.87%  "

$ws.Range("D15").Value = "'7.261"
$ws.Range("E15").Value = "  +1.22%  "

$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").Value = "'92.60"
$ws.Range("E17").Value = "  +0.73%  "

$ws.Range("E18").Value = "  +1.08%  "

$ws.Range("D19").Value = "'0.06514"
$ws.Range("E19").Value = "  +1.30%  "

$ws.Range("D20").Value = "'18.52"
$ws.Range("E20").Value = "  +3.39%  "

$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").Value = "'5.950"
$ws.Range("E22").Value = "  +2.50%  "

$ws.Range("D23").Value = "'30.156.16"
$ws.Range("E23").Value = "  +0.62%  "

$ws.Range("E24").Value = "  +2.47%  "

$ws.Range("D25").Value = "'2.194"
$ws.Range("E25").Value = "  +2.13%  "

$ws.Range("D26").Value = "'2.134.04"
$ws.Range("E26").Value = "  +2.57%  "

$ws.Range("D27").Value = "'21.90"
$ws.Range("E27").Value = "  +4.59%  "

$ws.Range("D28").Value = "'162.56"
$ws.Range("E28").Value = "  +1.00%  "

$ws.Range("D29").Value = "'2.268"
$ws.Range("E29").Value = "  +1.49%  "

$ws.Range("D30").Value = "'129.03"
$ws.Range("E30").Value = "  +1.60%  "

$ws.Range("D31").Value = "'1.130"
$ws.Range("E31").Value = "  +6.60%  "

$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("D33").Value = "'5.960"
$ws.Range("E33").Value = "  +0.97%  "

$ws.Range("D34").Value = "'3.791"
$ws.Range("E34").Value = "  +1.44%  "

$ws.Range("D35").Value = "'0.02455"
$ws.Range("E35").Value = "  +1.46%  "

$ws.Range("E36").Value = "  +1.81%  "

$ws.Range("D37").Value = "'0.06446"
$ws.Range("E37").Value = "  +1.58%  "

$ws.Range("D38").Value = "'0.2152"
$ws.Range("E38").Value = "  +0.36%  "

$ws.Range("D39").Value = "'0.6520"

$ws.Range("E40").Value = "  +2.41%  "

$ws.Range("D41").Value = "'8.593"
$ws.Range("E41").Value = "  +1.19%  "

$ws.Range("E42").Value = "  +1.96%  "

$ws.Range("E43").Value = "  +1.10%  "

$ws.Range("D44").Value = "'13.43"
$ws.Range("E44").Value = "  +3.72%  "

$ws.Range("D45").Value = "'2.196"
$ws.Range("E45").Value = "  +10.13%  "

$ws.Range("D46").Value = "'0.6058"
$ws.Range("E46").Value = "  +2.70%  "

$ws.Range("D47").Value = "'3.626"
$ws.Range("E47").Value = "  -0.22%  "

$ws.Range("D48").Value = "'1.211"
$ws.Range("E48").Value = "  +1.02%  "

$ws.Range("D49").Value = "'122.37"
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("E50").Value = "  +2.43%  "

$ws.Range("D51").Value = "'79.01"
$ws.Range("E51").Value = "  +3.08%  "
